$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data: Date 44930 (2023-01-04) and Time 2.
# Copy the date formatting from A5 onto A6 so it keeps the same cell
# style (m/d/yyyy number format + font), then set the new date value.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 44930

# B6 just needs the plain value - the column's default style already
# applies the correct formatting, so no paste-special is needed here
# (doing so would also break formula dependency tracking for B6).
$ws.Range("B6").Value = 2

# Update the D2 formula multiplier from 30 to 45 (set after the new
# row 6 values exist so SUM(B2:B100) recalculates including B6).
$ws.Range("D2").Formula = "=SUM(B2:B100)*45"

# Update the worksheet's active selection to D4, as in the diff.
$ws.Range("D4").Select()

$wb.Save()
